# Adding The Data Mining Final Revisions
#
# Slide 2 ("Weak 1 / Introduction") is duplicated to create a new third
# slide ("Tour 2 / Data Preprocessing") that reuses the same picture /
# decoration shapes, only the headline rectangle's position and text
# differ.

$p = $ppt.ActivePresentation

# Slide 2 is the template for the new slide.
$s2 = $p.Slides.Item(2)

# Duplicate() inserts the copy immediately after slide 2, i.e. as the new
# slide 3, and wires up the slide list / relationships automatically.
$s2.Duplicate()

$s3 = $p.Slides.Item(3)

# Locate the headline textbox ("Rectangle 1") on the new slide.
$headline = $s3.Shapes.Item("Rectangle 1")

# Reposition / resize it (left + width change; top + height stay the same).
$headline.Left = 46.15559387207031
$headline.Width = 447.6896850393701

# Update its two lines of text, keeping all existing run formatting.
$tf = $headline.TextFrame
$tr = $tf.TextRange
$tr.Paragraphs(1).Text = " Tour 2"
$tr.Paragraphs(2).Text = "Data Preprocessing"
